$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D25").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E25").Value = '[''Normal'']'

$ws.Range("D29").Value = '[0, 0, 0, 1, 0, 0, 1]'
$ws.Range("E29").Value = '[''ParamViolation'', ''SoftwareFault'']'

$ws.Range("D39").Value = '[1, 0, 1, 0, 0, 0, 1]'
$ws.Range("E39").Value = '[''Normal'', ''HardwareFault'', ''SoftwareFault'']'

$ws.Range("D41").Value = '[1, 0, 1, 0, 0, 0, 0]'
$ws.Range("E41").Value = '[''Normal'', ''HardwareFault'']'

$ws.Range("D53").Value = '[1, 0, 1, 0, 0, 0, 0]'
$ws.Range("E53").Value = '[''Normal'', ''HardwareFault'']'

$ws.Range("D54").Value = '[0, 0, 0, 0, 0, 1, 0]'
$ws.Range("E54").Value = '[''CommunicationIssue'']'

$ws.Range("D68").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E68").Value = '[''Normal'']'

$ws.Range("D69").Value = '[1, 1, 0, 0, 0, 0, 0]'
$ws.Range("E69").Value = '[''Normal'', ''SurroundingEnvironment'']'

$ws.Range("D74").Value = '[1, 0, 0, 0, 0, 0, 1]'
$ws.Range("E74").Value = '[''Normal'', ''SoftwareFault'']'

$ws.Range("D83").Value = '[1, 1, 0, 0, 0, 0, 0]'
$ws.Range("E83").Value = '[''Normal'', ''SurroundingEnvironment'']'

$ws.Range("D84").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E84").Value = '[''Normal'']'

$ws.Range("D92").Value = '[1, 0, 1, 0, 0, 0, 1]'
$ws.Range("E92").Value = '[''Normal'', ''HardwareFault'', ''SoftwareFault'']'

$ws.Range("D107").Value = '[1, 0, 0, 0, 0, 1, 0]'
$ws.Range("E107").Value = '[''Normal'', ''CommunicationIssue'']'

$ws.Range("D109").Value = '[1, 1, 0, 0, 0, 0, 0]'
$ws.Range("E109").Value = '[''Normal'', ''SurroundingEnvironment'']'

$ws.Range("D113").Value = '[1, 0, 1, 0, 0, 0, 0]'
$ws.Range("E113").Value = '[''Normal'', ''HardwareFault'']'
